$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.471.31'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.249.60'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '564.88'
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.59'
$ws.Range('E6').Value = '  -4.75%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.579'
$ws.Range('E8').Value = '  +1.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.240.05'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.564'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.84'
$ws.Range('E12').Value = '  -3.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000265'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '693.48'
$ws.Range('E14').Value = '  +11.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.774.18'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.22'
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.671.33'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.264.37'
$ws.Range('E19').Value = '  -0.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.14'
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.60'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.877'
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.75'
$ws.Range('E23').Value = '  -6.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.06'
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.99'
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.83'
$ws.Range('E26').Value = '  -2.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.67'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.16'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '32.42'
$ws.Range('E29').Value = '  +6.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.30'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.60'
$ws.Range('E31').Value = '  +2.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '572.06'
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.823.69'
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.66'
$ws.Range('E34').Value = '  -1.45%  '
$ws.Range('B35').Value = 'Dai'
$ws.Range('C35').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.102'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.08'
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.23'
$ws.Range('E38').Value = '  -11.41%  '
$ws.Range('E39').Value = '  +1.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.33'
$ws.Range('E40').Value = '  -2.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.56'
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '31.35'
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.04'
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0658'
$ws.Range('E44').Value = '  -2.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.323'
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0401'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.36'
$ws.Range('E49').Value = '  +8.04%  '
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '128.16'
$ws.Range('E51').Value = '  -0.71%  '
